# Updated cryptos list with GitHub Actions
# Applies the per-cell value changes described by the diff against
# the "cryptos" worksheet (Coin / Link / Price / Volume(1h) columns).
# Numeric-looking Price strings are written with a leading apostrophe and
# the cell style is reset to "Normal" afterwards so Excel keeps them as
# plain text (matching the source `inlineStr` cells) instead of silently
# converting them to floating point numbers, while not leaving any stray
# number-format/style behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = '43.136.92'
$ws.Range("E2").Value = '  -4.90%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '2.229.09'
$ws.Range("E3").Value = '  -5.87%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  +0.14%  '

# Row 5 - BNB
$ws.Range("D5").Value = '''315.92'
$ws.Range("D5").Style = "Normal"

# Row 6 - Solana
$ws.Range("D6").Value = '''100.82'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.37%  '

# Row 7 - XRP
$ws.Range("D7").Value = '''0.585'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -6.96%  '

# Row 8 - USDC
$ws.Range("E8").Value = '  -0.08%  '

# Row 9 - Cardano
$ws.Range("D9").Value = '''0.560'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -8.05%  '

# Row 10 - Avalanche
$ws.Range("D10").Value = '''36.84'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -9.48%  '

# Row 11 - OKB
$ws.Range("D11").Value = '''54.44'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.98%  '

# Row 12 - Dogecoin
$ws.Range("E12").Value = '  -10.26%  '

# Row 13 - Polkadot
$ws.Range("D13").Value = '''7.65'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -9.49%  '

# Row 14 - TRON
$ws.Range("E14").Value = '  -1.40%  '

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = '2.569.49'
$ws.Range("E15").Value = '  -5.83%  '

# Row 16 - Polygon
$ws.Range("D16").Value = '''0.858'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -12.24%  '

# Row 17 - Chainlink
$ws.Range("D17").Value = '''14.32'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -6.43%  '

# Row 18 - WrappedEther
$ws.Range("D18").Value = '2.227.44'
$ws.Range("E18").Value = '  -5.79%  '

# Row 19 - WrappedBTC
$ws.Range("D19").Value = '43.086.97'
$ws.Range("E19").Value = '  -5.07%  '

# Row 20 - InternetComputer(DFINITY)
$ws.Range("D20").Value = '''14.26'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.69%  '

# Row 21 - ShibaInu
$ws.Range("E21").Value = '  -9.67%  '

# Row 22 - Uniswap
$ws.Range("D22").Value = '''6.48'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -10.07%  '

# Row 23 - Litecoin
$ws.Range("D23").Value = '''65.39'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -10.67%  '

# Row 24 - PancakeSwap
$ws.Range("E24").Value = '  -9.74%  '

# Row 25 - BitcoinCash
$ws.Range("D25").Value = '''237.55'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -8.54%  '

# Row 26 - ImmutableX
$ws.Range("E26").Value = '  -9.93%  '

# Row 27 - Dai
$ws.Range("E27").Value = '  +0.06%  '

# Row 28 - LEO
$ws.Range("D28").Value = '''4.07'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.25%  '

# Row 29 - Cosmos->Toncoin
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '''2.23'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.55%  '

# Row 30 - Toncoin->Cosmos
$ws.Range("B30").Value = 'Cosmos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D30").Value = '''9.98'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -9.65%  '

# Row 31 - Filecoin
$ws.Range("D31").Value = '''6.40'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -12.00%  '

# Row 32 - EthereumClassic
$ws.Range("D32").Value = '''20.45'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -8.17%  '

# Row 33 - InjectiveProtocol
$ws.Range("D33").Value = '''34.40'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.41%  '

# Row 34 - Hedera
$ws.Range("D34").Value = '''0.0870'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -10.29%  '

# Row 35 - Monero
$ws.Range("D35").Value = '''154.20'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.37%  '

# Row 36 - WEMIXToken
$ws.Range("D36").Value = '''2.77'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.89%  '

# Row 37 - LidoDAOToken
$ws.Range("D37").Value = '''3.18'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +8.54%  '

# Row 38 - ARBITRUM
$ws.Range("E38").Value = '  +2.36%  '

# Row 39 - Stellar
$ws.Range("D39").Value = '''0.121'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.78%  '

# Row 40 - RenderToken
$ws.Range("D40").Value = '''4.42'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.29%  '

# Row 41 - Kaspa
$ws.Range("E41").Value = '  -11.56%  '

# Row 42 - NEARProtocol
$ws.Range("D42").Value = '''3.69'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.94%  '

# Row 43 - VeChain
$ws.Range("E43").Value = '  -8.64%  '

# Row 44 - Celestia
$ws.Range("D44").Value = '''13.04'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.30%  '

# Row 45 - FirstDigitalUSD
$ws.Range("E45").Value = '  -0.02%  '

# Row 46 - Maker
$ws.Range("D46").Value = '1.797.85'
$ws.Range("E46").Value = '  -0.71%  '

# Row 47 - BitcoinSV
$ws.Range("D47").Value = '''87.99'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -10.41%  '

# Row 48 - Algorand
$ws.Range("E48").Value = '  -9.19%  '

# Row 49 - ordi
$ws.Range("D49").Value = '''76.83'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -8.00%  '

# Row 50 - THORChain
$ws.Range("D50").Value = '''5.32'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -8.61%  '

# Row 51 - MultiversX
$ws.Range("D51").Value = '''59.00'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -15.74%  '
